# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets now that the
# handback has completed, updates the Overview sheet's status text, and
# widens a few columns so the new/longer values are readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet: status text moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both locale status columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$handedBack = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# Widen the zh-cn / de-de status columns to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# 2) zh-cn sheet: fill in target/handback file + handback datetime for
#    both rows, and turn the target-file cells into hyperlinks (matching
#    the existing "Source File Name" hyperlink style/behaviour).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$target59199 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/59199e9b-216c-4568-9c85-c61bf9ca802a.md"
$targetCc515 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $target59199, [Type]::Missing, [Type]::Missing, "59199e9b-216c-4568-9c85-c61bf9ca802a.md")
$wsZhCn.Range("J2").Value = "59199e9b-216c-4568-9c85-c61bf9ca802a.bb4f61a0eeaf7b4a1715626c371b57790a83e8c1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 23:05:30"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetCc515, [Type]::Missing, [Type]::Missing, "cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md")
$wsZhCn.Range("J3").Value = "cc515795-00a0-4335-9ec6-ff4d30aa0e5e.25ed7a593b7068aebc82aa7099dcafe3e64582a2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-20 23:05:30"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# 3) de-de sheet: same shape as zh-cn, but with its own handback
#    datetime (distinct from zh-cn's).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $target59199, [Type]::Missing, [Type]::Missing, "59199e9b-216c-4568-9c85-c61bf9ca802a.md")
$wsDeDe.Range("J2").Value = "59199e9b-216c-4568-9c85-c61bf9ca802a.bb4f61a0eeaf7b4a1715626c371b57790a83e8c1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 23:05:37"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetCc515, [Type]::Missing, [Type]::Missing, "cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md")
$wsDeDe.Range("J3").Value = "cc515795-00a0-4335-9ec6-ff4d30aa0e5e.25ed7a593b7068aebc82aa7099dcafe3e64582a2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-20 23:05:37"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667
